# gmail_flight.xlsx update
# - Fix mislabeled airline text in D12, D13, D19 ("EVA cooperated")
# - Append 10 new tracked-price rows (20-29) pulled from the latest snippets
# - Extend the sheet selection to cover the new data range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing rows: airline column only (plain text, no reformat needed) ---
$ws.Range("D12").Value = "EVA cooperated"
$ws.Range("D13").Value = "EVA cooperated"
$ws.Range("D19").Value = "EVA cooperated"

# --- Append new rows 20-29 ---
# Pre-format as Text so date-like / currency-like strings are kept verbatim
# (e.g. "09/28/19" and "$1,301") instead of being parsed into numbers.
$ws.Range("A20:D29").NumberFormat = "@"

$newRows = @(
  @("09/28/19", '$1,301', '$1,147', "EVA cooperated"),
  @("09/29/19", '$742',   '$848',   "Multiple airlines"),
  @("09/30/19", '$1,147', '$1,301', "EVA cooperated"),
  @("10/01/19", '$878',   '$742',   "EVA Air & cooeperated"),
  @("10/04/19", '$825',   '$878',   "Multiple airlines"),
  @("10/05/19", '$1,353', '$2,028', "EVA Air"),
  @("10/06/19", '$1,966', '$1,353', "EVA Air"),
  @("10/07/19", '$2,520', '$1,966', "EVA Air"),
  @("10/08/19", '$2,365', '$2,520', "EVA Air"),
  @("10/09/19", '$757',   '$825',   "Multiple airlines")
)

$r = 20
foreach ($row in $newRows) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $r++
}

# --- Update the visible selection to match the new used range ---
$ws.Activate() | Out-Null
$ws.Range("A2:D29").Select() | Out-Null
